$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set the Runmode column (C) to "N" for TestCase_B2 through TestCase_B6 (rows 3-7),
# leaving TestCase_B1 (row 2) as "Y" so only A1/B1 keep running.
$ws.Range("C3:C7").Value = "N"

# Move the active selection to C8 (a single cell, no range) to reflect the
# cursor position after the edit.
$ws.Range("C8").Select()

$wb.Save()
